$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the hyperlinks that live on row 4 (C4/D4) before the row itself
# disappears, so the saved package doesn't retain orphaned <hyperlink>
# entries / relationship targets for the removed reviewer emails.
$changed = $true
while ($changed) {
    $changed = $false
    foreach ($hl in $ws.Hyperlinks) {
        $addr = $hl.Range.Address()
        if ($addr -eq '$C$4' -or $addr -eq '$D$4') {
            $hl.Delete()
            $changed = $true
            break
        }
    }
}

# Delete row 4 entirely, shifting rows 5-7 up to 4-6. This drops the
# review record that lived there and, because the only remaining
# references to those shared strings disappear with it, they get pruned
# from sharedStrings.xml on save.
$ws.Rows("4").Delete()

$ws.Range("A4").Select()
